# Countries table edit: the "Data" sheet's D6 cell had a trailing "; " left
# over from data entry; trim it down to match the other entries in the
# column and leave the cell selected (mirrors the author's manual edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "usa; uk; australia; india; germany; spain; france; russia; china; new zealand; canada"

$ws.Range("D6").Select()
